$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: restyle row 15 to the "standalone entry" box-border style (same as rows 6-8) ---
$ws.Range("A6:E6").Copy() | Out-Null
$ws.Range("A15:E15").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Step 2: add row 16 (first row of a 2-row entry) using the same format as row 9 / row 11 ---
$ws.Range("A11:E11").Copy() | Out-Null
$ws.Range("A16:E16").PasteSpecial(-4122) | Out-Null

# --- Step 3: add row 17 (closing/filename-only row of the entry) using the same format as row 10 / row 12 ---
$ws.Range("A12:E12").Copy() | Out-Null
$ws.Range("A17:E17").PasteSpecial(-4122) | Out-Null

# --- Step 4: add row 18 (standalone entry) using the same format as row 6-8 ---
$ws.Range("A6:E6").Copy() | Out-Null
$ws.Range("A18:E18").PasteSpecial(-4122) | Out-Null

# --- Step 5: add row 19 (last, standalone/no-border-close entry) using the same format as row 9 / row 11 ---
$ws.Range("A11:E11").Copy() | Out-Null
$ws.Range("A19:E19").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Step 6: populate cell values, in the exact order the new shared strings appear in the target file ---
$ws.Range("C16").Value = " I enjoy getting [CS:I]Prize Ticket[CR]" + [char]0x2026
$ws.Range("A16").Value = " SCRIPT/P01P04A/us0403.ssb"
$ws.Range("D16").Value = " Люблю получать [CS:I]Призовые Билеты[CR]..."
$ws.Range("E16").Value = " Ìýáìý ðïìôœàóû [CS:I]Ðñéèïâúå Áéìåóú[CR]..."
$ws.Range("B16").Value = 57

$ws.Range("A17").Value = "SCRIPT/P01P04A/us2008.ssb"

$ws.Range("C18").Value = " [CS:N]Drowzee[CR]?[K] ...Oh, right.\nThat petty crook."
$ws.Range("A18").Value = "SCRIPT/T01P02A/us2010.ssb"
$ws.Range("D18").Value = " [CS:N]Дроузи[CR]?[K] ...Ах, точно. Тот\nжалкий негодяй."
$ws.Range("E18").Value = " [CS:N]Äñïôèé[CR]?[K] ...Àö, óïœîï. Óïó\nçàìëéê îåãïäÿê."
$ws.Range("B18").Value = 38

$ws.Range("A19").Value = "SCRIPT/T01P02A/us2014.ssb"
$ws.Range("C19").Value = " Why'd you bring in a known\noutlaw like [CS:N]Drowzee[CR]?"
$ws.Range("D19").Value = " Зачем вы привели сюда печально\nизвестного негодяя [CS:N]Дроузи[CR]?"
$ws.Range("E19").Value = " Èàœåí âú ðñéâåìé òýäà ðåœàìûîï\néèâåòóîïãï îåãïäÿÿ [CS:N]Äñïôèé[CR]?"
$ws.Range("B19").Value = 18

# --- Step 7: row heights to match the authored worksheet ---
$ws.Rows.Item(16).RowHeight = 57.6
$ws.Rows.Item(17).RowHeight = 43.2
$ws.Rows.Item(18).RowHeight = 43.2
$ws.Rows.Item(19).RowHeight = 43.2

# --- Step 8: view state - scroll position and active selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D19").Select() | Out-Null

Write-Host "edit complete"
